$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update spring hardpoint values (new MF-Swift v2212 hardpoints) ---
$ws.Range("G5").Value = 0.62
$ws.Range("H5").Value = 0.65

$ws.Range("G6").Value = 0.85
$ws.Range("H6").Value = 0.19

$ws.Range("G9").Value = 0.62
$ws.Range("H9").Value = 0.65

$ws.Range("G10").Value = 0.85
$ws.Range("H10").Value = 0.19

# --- Re-format the x/y/z columns to 2 decimal places (style 19) ---
# F9/F10 change format only (values unchanged); G/H columns already
# rewritten above, now align all of them on the new "0.00" number format.
# (Set one cell at a time -- a multi-area Range.NumberFormat assignment
# only reliably applies to the first area in this host.)
$ws.Range("F9").NumberFormat = "0.00"
$ws.Range("F10").NumberFormat = "0.00"
$ws.Range("G5").NumberFormat = "0.00"
$ws.Range("H5").NumberFormat = "0.00"
$ws.Range("G6").NumberFormat = "0.00"
$ws.Range("H6").NumberFormat = "0.00"
$ws.Range("G9").NumberFormat = "0.00"
$ws.Range("H9").NumberFormat = "0.00"
$ws.Range("G10").NumberFormat = "0.00"
$ws.Range("H10").NumberFormat = "0.00"
